$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: prepare row 13 as a duplicate of row 12 (copies formatting: date/time/text styles) ---
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)

# Put a temporary copy of C12's current text into C13 so both cells share the
# same underlying string before we diverge them (needed so the shared-string
# table ends up with the modified text at a new slot, and the brand new text
# re-using the freed slot).
$ws.Range("C13").Value = $ws.Range("C12").Value

# --- Step 2: diverge the two cells ---
# C13 becomes the brand-new journal entry text first...
$ws.Range("C13").Value = "M. Benzonana m'a aidé à corriger les bugs qui m'empèchait d'avoir un environnement de dévelloppement fonctionnel"

# ...then C12 is updated in place to the corrected wording of the original entry.
$ws.Range("C12").Value = "Préparation de l'environement de travail à la maison (sans WAMP) avec l'aide de Niclass Dorian"

# --- Step 3: fill in the rest of row 13 ---
$ws.Range("A13").Value = "03/19/2020"
$ws.Range("B13").Value = 5.25
$ws.Rows.Item(13).RowHeight = 30

# --- Step 4: add row 14 (new entry) ---
$ws.Range("A12:C12").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Value = "03/20/2020"
$ws.Range("B14").Value = 5.25
$ws.Range("C14").Value = "M. Benzonana et Meylan Benoit m'ont aidé à appliquer la fonctionnalité des calculs aléatoires"
$ws.Rows.Item(14).RowHeight = 30

# --- Step 5: add row 15 (new entry, whole-number time like row 2/8 => style s=5) ---
$ws.Range("A12:C12").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A15").Value = "03/24/2020"
$ws.Range("B2:B2").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "Réorganisation personnel sur le projet et avancement de la documentation"
$ws.Rows.Item(15).RowHeight = 30

# --- Step 6: selection / dimension bookkeeping ---
$ws.Range("C15").Select()
